$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column R (overload values) for rows 5-21 with a 0.6 .. 2.2 ramp (step 0.1)
$values = @(0.6, 0.7, 0.8, 0.9, 1, 1.1, 1.2, 1.3, 1.4, 1.5, 1.6, 1.7, 1.8, 1.9, 2, 2.1, 2.2)
$row = 5
foreach ($v in $values) {
    $ws.Cells.Item($row, 18).Value = $v
    $row++
}

# Update the view state: top-left visible cell and active selection
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("R24").Select()
